# "DAC buf disabled" / "Current setup and measure constants adjusted"
#
# Sheet1 was blank; populate it with the "Current setup" table:
#   row1  Current setup                    (Accent1 banner)
#   row2  Dac1      500
#   row3  Ur100,mV  53.2
#   row4  I1,uA     =B3*10
#   row5  DAC2      3000
#   row6  Ur100,mV  298.39999999999998
#   row7  I2,uA     =B6*10
#   row8  Divider   =B7-B4
#   row9  a         =B5-B2
#   row10 b         =B2*B8-B4*B9
#
# Cell values/styles are written in the same order the original author
# would have typed them (header -> Dac1 block -> DAC2 block -> derived
# rows, with the repeated "Ur100,mV" label typed in last) so that the
# shared-string table and the style gallery end up built in the same
# order as the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A a bit wider than default so the labels fit.
$ws.Columns.Item(1).ColumnWidth = 10.2

# Row 1 - section banner
$ws.Range("A1").Value = "Current setup"
$ws.Range("A1").Style = "Accent1"

# Row 2 - Dac1 setting
$ws.Range("A2").Value = "Dac1"
$ws.Range("B2").Value = 500
$ws.Range("A2:B2").Style = "Good"

# Row 3 - measured Ur100 (value only for now, label typed later)
$ws.Range("B3").Value = 53.2
$ws.Range("A3:B3").Style = "Good"

# Row 4 - computed I1
$ws.Range("A4").Value = "I1,uA"
$ws.Range("B4").Formula = "=B3*10"
$ws.Range("A4:B4").Style = "Neutral"

# Row 5 - DAC2 setting
$ws.Range("A5").Value = "DAC2"
$ws.Range("B5").Value = 3000
$ws.Range("A5:B5").Style = "Good"

# Row 6 - measured Ur100 (second occurrence)
$ws.Range("B6").Value = 298.39999999999998
$ws.Range("A6:B6").Style = "Good"

# Row 7 - computed I2, underlined to separate the raw setup from the
# derived divider/line constants below.
$ws.Range("A7").Value = "I2,uA"
$ws.Range("B7").Formula = "=B6*10"
$ws.Range("A7:B7").Style = "Neutral"
$ws.Range("A7:B7").Borders.Item(9).LineStyle = 1

# Row 8 - Divider
$ws.Range("A8").Value = "Divider"
$ws.Range("B8").Formula = "=B7-B4"
$ws.Range("A8:B8").Style = "Calculation"

# Row 9 - a
$ws.Range("A9").Value = "a"
$ws.Range("B9").Formula = "=B5-B2"
$ws.Range("A9:B9").Style = "Calculation"

# Row 10 - b
$ws.Range("A10").Value = "b"
$ws.Range("B10").Formula = "=B2*B8-B4*B9"
$ws.Range("A10:B10").Style = "Calculation"

# The repeated "Ur100,mV" label typed last, re-using the shared string
# created here for both occurrences.
$ws.Range("A3").Value = "Ur100,mV"
$ws.Range("A6").Value = "Ur100,mV"

# Final selection / page setup, matching the saved workbook state.
$ws.Range("B7").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
